$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44978
$ws.Range("J2").Value2 = 1000
$ws.Range("K2").Value2 = 1800
$ws.Range("L2").Value2 = 2000
$ws.Range("M2").Value2 = 1900
$ws.Range("P2").Value2 = 633

$ws.Range("D3").Value2 = 44965
$ws.Range("J3").Value2 = 1120
$ws.Range("K3").Value2 = 2000
$ws.Range("L3").Value2 = 2500
$ws.Range("M3").Value2 = 2250
$ws.Range("P3").Value2 = 750

$ws.Range("D4").Value2 = 45006
$ws.Range("J4").Value2 = 1100
$ws.Range("K4").Value2 = 2000
$ws.Range("L4").Value2 = 2500
$ws.Range("M4").Value2 = 2250
$ws.Range("P4").Value2 = 750

$ws.Range("D5").Value2 = 44848
$ws.Range("J5").Value2 = 1000
$ws.Range("K5").Value2 = 1500
$ws.Range("L5").Value2 = 2000
$ws.Range("M5").Value2 = 1750
$ws.Range("P5").Value2 = 583

$ws.Range("D6").Value2 = 45020
$ws.Range("J6").Value2 = 1200
$ws.Range("K6").Value2 = 2000
$ws.Range("L6").Value2 = 2500
$ws.Range("M6").Value2 = 2250
$ws.Range("P6").Value2 = 750

$ws.Range("D7").Value2 = 44985
$ws.Range("J7").Value2 = 1000
$ws.Range("K7").Value2 = 2000
$ws.Range("L7").Value2 = 2500
$ws.Range("M7").Value2 = 2250
$ws.Range("P7").Value2 = 750

$ws.Range("D9").Value2 = 45070
$ws.Range("J9").Value2 = 800
$ws.Range("K9").Value2 = 2000
$ws.Range("L9").Value2 = 2500
$ws.Range("M9").Value2 = 2250
$ws.Range("P9").Value2 = 750

$ws.Range("D10").Value2 = 44999
$ws.Range("J10").Value2 = 1100
$ws.Range("K10").Value2 = 2000
$ws.Range("L10").Value2 = 2500
$ws.Range("M10").Value2 = 2250
$ws.Range("P10").Value2 = 750

$ws.Range("D11").Value2 = 44971
$ws.Range("J11").Value2 = 1000
$ws.Range("K11").Value2 = 2000
$ws.Range("L11").Value2 = 2500
$ws.Range("M11").Value2 = 2250
$ws.Range("P11").Value2 = 750

$ws.Range("D12").Value2 = 44992
$ws.Range("J12").Value2 = 1040
$ws.Range("K12").Value2 = 2000
$ws.Range("L12").Value2 = 2500
$ws.Range("M12").Value2 = 2250
$ws.Range("P12").Value2 = 750

$ws.Range("D13").Value2 = 45062
$ws.Range("J13").Value2 = 1100
$ws.Range("K13").Value2 = 2000
$ws.Range("L13").Value2 = 2500
$ws.Range("M13").Value2 = 2250
$ws.Range("P13").Value2 = 750

$ws.Range("D14").Value2 = 44827
$ws.Range("J14").Value2 = 1200
$ws.Range("K14").Value2 = 2000
$ws.Range("L14").Value2 = 2500
$ws.Range("M14").Value2 = 2250
$ws.Range("P14").Value2 = 750

$ws.Range("D15").Value2 = 44911
$ws.Range("J15").Value2 = 700
$ws.Range("K15").Value2 = 1800
$ws.Range("L15").Value2 = 2000
$ws.Range("M15").Value2 = 1900
$ws.Range("P15").Value2 = 633

$ws.Range("D16").Value2 = 44970
$ws.Range("J16").Value2 = 800
$ws.Range("K16").Value2 = 2000
$ws.Range("L16").Value2 = 2500
$ws.Range("M16").Value2 = 2250
$ws.Range("P16").Value2 = 750

$ws.Range("D17").Value2 = 45035
$ws.Range("J17").Value2 = 1100
$ws.Range("K17").Value2 = 2000
$ws.Range("L17").Value2 = 2500
$ws.Range("M17").Value2 = 2250
$ws.Range("P17").Value2 = 750

$ws.Range("D18").Value2 = 45034
$ws.Range("J18").Value2 = 1100
$ws.Range("K18").Value2 = 2000
$ws.Range("L18").Value2 = 2500
$ws.Range("M18").Value2 = 2250
$ws.Range("P18").Value2 = 750

$ws.Range("D19").Value2 = 44964
$ws.Range("J19").Value2 = 1000
$ws.Range("K19").Value2 = 2000
$ws.Range("L19").Value2 = 2500
$ws.Range("M19").Value2 = 2250
$ws.Range("P19").Value2 = 750

$ws.Range("D20").Value2 = 45041
$ws.Range("J20").Value2 = 1160
$ws.Range("K20").Value2 = 2000
$ws.Range("L20").Value2 = 2500
$ws.Range("M20").Value2 = 2250
$ws.Range("P20").Value2 = 750

$ws.Range("D21").Value2 = 45028
$ws.Range("J21").Value2 = 1000
$ws.Range("K21").Value2 = 2000
$ws.Range("L21").Value2 = 2500
$ws.Range("M21").Value2 = 2250
$ws.Range("P21").Value2 = 750

$ws.Range("D22").Value2 = 45013
$ws.Range("J22").Value2 = 1100
$ws.Range("K22").Value2 = 2000
$ws.Range("L22").Value2 = 2500
$ws.Range("M22").Value2 = 2250
$ws.Range("P22").Value2 = 750

$ws.Range("D23").Value2 = 44685
$ws.Range("J23").Value2 = 400
$ws.Range("K23").Value2 = 1500
$ws.Range("L23").Value2 = 2000
$ws.Range("M23").Value2 = 1750
$ws.Range("P23").Value2 = 583

$ws.Range("D24").Value2 = 45084
$ws.Range("J24").Value2 = 900
$ws.Range("K24").Value2 = 2000
$ws.Range("L24").Value2 = 2500
$ws.Range("M24").Value2 = 2250
$ws.Range("P24").Value2 = 750

$ws.Range("D25").Value2 = 44883
$ws.Range("J25").Value2 = 500
$ws.Range("K25").Value2 = 1800
$ws.Range("L25").Value2 = 2000
$ws.Range("M25").Value2 = 1900
$ws.Range("P25").Value2 = 633

$ws.Range("D26").Value2 = 44951
$ws.Range("J26").Value2 = 800
$ws.Range("K26").Value2 = 2000
$ws.Range("L26").Value2 = 2500
$ws.Range("M26").Value2 = 2250
$ws.Range("P26").Value2 = 750

$ws.Range("D27").Value2 = 45091
$ws.Range("J27").Value2 = 800
$ws.Range("K27").Value2 = 2000
$ws.Range("L27").Value2 = 2500
$ws.Range("M27").Value2 = 2250
$ws.Range("P27").Value2 = 750

$ws.Range("D28").Value2 = 45007
$ws.Range("J28").Value2 = 1160
$ws.Range("K28").Value2 = 2000
$ws.Range("L28").Value2 = 2500
$ws.Range("M28").Value2 = 2250
$ws.Range("P28").Value2 = 750

$ws.Range("D29").Value2 = 44910
$ws.Range("J29").Value2 = 1000
$ws.Range("K29").Value2 = 1800
$ws.Range("L29").Value2 = 2000
$ws.Range("M29").Value2 = 1900
$ws.Range("P29").Value2 = 633

$ws.Range("D30").Value2 = 45077
$ws.Range("J30").Value2 = 760
$ws.Range("K30").Value2 = 2000
$ws.Range("L30").Value2 = 2500
$ws.Range("M30").Value2 = 2250
$ws.Range("P30").Value2 = 750

$ws.Range("D31").Value2 = 44953
$ws.Range("J31").Value2 = 1000
$ws.Range("K31").Value2 = 2000
$ws.Range("L31").Value2 = 2500
$ws.Range("M31").Value2 = 2250
$ws.Range("P31").Value2 = 750
